$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 282.72223
$ws.Range("I33").Value = 286.875
$ws.Range("K33").Value = 286.875
$ws.Range("M33").Value = -57.875
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H76").Value = 3316.6667
$ws.Range("I76").Value = 3321.4285
$ws.Range("J76").Value = 3300
$ws.Range("K76").Value = 3321.4285
$ws.Range("L76").Value = 3300
$ws.Range("M76").Value = -3006.4285
$ws.Range("N76").Value = -3930
$ws.Range("H79").Value = 3316.6667
$ws.Range("I79").Value = 3321.4285
$ws.Range("J79").Value = 3300
$ws.Range("K79").Value = 3321.4285
$ws.Range("L79").Value = 3300
$ws.Range("M79").Value = -2229.4285
$ws.Range("N79").Value = -5484
$ws.Range("H116").Value = 2828.0557
$ws.Range("I116").Value = 2814.6428
$ws.Range("J116").Value = 2875
$ws.Range("K116").Value = 2814.6428
$ws.Range("L116").Value = 2875
$ws.Range("M116").Value = 627.3571999999999
$ws.Range("N116").Value = -9759
$ws.Range("H137").Value = 1750.8055
$ws.Range("I137").Value = 1773.7646
$ws.Range("K137").Value = 5321.293799999999
$ws.Range("M137").Value = -2771.293799999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5328618.5
$ws.Range("I32").Value = 6105385
$ws.Range("J32").Value = 20713.416
$ws.Range("K32").Value = 6105385
$ws.Range("L32").Value = 20713.416
$ws.Range("M32").Value = -6105098
$ws.Range("N32").Value = -21287.416
$ws.Range("H43").Value = 11657.125
$ws.Range("I43").Value = 3500
$ws.Range("J43").Value = 12822.429
$ws.Range("K43").Value = 3500
$ws.Range("L43").Value = 12822.429
$ws.Range("M43").Value = -3187
$ws.Range("N43").Value = -13448.429
$ws.Range("H74").Value = 2409.425
$ws.Range("I74").Value = 1696.0416
$ws.Range("J74").Value = 3479.5
$ws.Range("K74").Value = 1696.0416
$ws.Range("L74").Value = 3479.5
$ws.Range("M74").Value = -822.0416
$ws.Range("N74").Value = -5227.5
$ws.Range("H77").Value = 2409.425
$ws.Range("I77").Value = 1696.0416
$ws.Range("J77").Value = 3479.5
$ws.Range("K77").Value = 8480.208000000001
$ws.Range("L77").Value = 17397.5
$ws.Range("M77").Value = -4112.208000000001
$ws.Range("N77").Value = -26133.5
$ws.Range("H103").Value = 57785.668
$ws.Range("J103").Value = 57785.668
$ws.Range("L103").Value = 57785.668
$ws.Range("N103").Value = -60129.668

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 329.83334
$ws.Range("I94").Value = 235.8
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 235.8
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 215.2
$ws.Range("N94").Value = -1702

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.53846
$ws.Range("I7").Value = 52.25
$ws.Range("J7").Value = 73.8
$ws.Range("K7").Value = 52.25
$ws.Range("L7").Value = 73.8
$ws.Range("M7").Value = 60.75
$ws.Range("N7").Value = -299.8
$ws.Range("H31").Value = 6131.1846
$ws.Range("I31").Value = 1843.2858
$ws.Range("J31").Value = 8177.6816
$ws.Range("K31").Value = 1843.2858
$ws.Range("L31").Value = 8177.6816
$ws.Range("M31").Value = -1548.2858
$ws.Range("N31").Value = -8767.6816
$ws.Range("H34").Value = 6131.1846
$ws.Range("I34").Value = 1843.2858
$ws.Range("J34").Value = 8177.6816
$ws.Range("K34").Value = 1843.2858
$ws.Range("L34").Value = 8177.6816
$ws.Range("M34").Value = -1641.2858
$ws.Range("N34").Value = -8581.6816
$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 332.57144
$ws.Range("I47").Value = 167.6
$ws.Range("J47").Value = 745
$ws.Range("K47").Value = 502.8
$ws.Range("L47").Value = 2235
$ws.Range("M47").Value = -71.79999999999995
$ws.Range("N47").Value = -3097
$ws.Range("H74").Value = 1013
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1013
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 2941.7021
$ws.Range("I122").Value = 463.22223
$ws.Range("K122").Value = 4169.00007
$ws.Range("M122").Value = -1719.00007

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1881800
$ws.Range("I80").Value = 2252250
$ws.Range("J80").Value = 400000
$ws.Range("K80").Value = 2252250
$ws.Range("L80").Value = 400000
$ws.Range("M80").Value = -2251252
$ws.Range("N80").Value = -401996
$ws.Range("H83").Value = 1881800
$ws.Range("I83").Value = 2252250
$ws.Range("J83").Value = 400000
$ws.Range("K83").Value = 11261250
$ws.Range("L83").Value = 2000000
$ws.Range("M83").Value = -11256258
$ws.Range("N83").Value = -2009984
$ws.Range("H114").Value = 39446.547
$ws.Range("J114").Value = 39446.547
$ws.Range("L114").Value = 39446.547
$ws.Range("N114").Value = -48124.547
$ws.Range("H132").Value = 2171.4285
$ws.Range("I132").Value = 1708.3334
$ws.Range("J132").Value = 3453.8462
$ws.Range("K132").Value = 5125.0002
$ws.Range("L132").Value = 10361.5386
$ws.Range("M132").Value = -2595.0002
$ws.Range("N132").Value = -15421.5386

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10376
$ws.Range("H82").Value = 3027.2727
$ws.Range("I82").Value = 4666.6665
$ws.Range("J82").Value = 2412.5
$ws.Range("K82").Value = 4666.6665
$ws.Range("L82").Value = 2412.5
$ws.Range("M82").Value = -4305.6665
$ws.Range("N82").Value = -3134.5
$ws.Range("H85").Value = 3027.2727
$ws.Range("I85").Value = 4666.6665
$ws.Range("J85").Value = 2412.5
$ws.Range("K85").Value = 4666.6665
$ws.Range("L85").Value = 2412.5
$ws.Range("M85").Value = -3418.6665
$ws.Range("N85").Value = -4908.5
$ws.Range("H93").Value = 10831.833
$ws.Range("I93").Value = 11362
$ws.Range("K93").Value = 11362
$ws.Range("M93").Value = -10114
$ws.Range("H103").Value = 98602
$ws.Range("J103").Value = 98602
$ws.Range("L103").Value = 98602
$ws.Range("N103").Value = -100946
$ws.Range("H106").Value = 79790
$ws.Range("J106").Value = 79790
$ws.Range("L106").Value = 79790
$ws.Range("N106").Value = -82314
$ws.Range("H132").Value = 3189.3914
$ws.Range("I132").Value = 2776
$ws.Range("J132").Value = 4677.6
$ws.Range("K132").Value = 8328
$ws.Range("L132").Value = 14032.8
$ws.Range("M132").Value = -5798
$ws.Range("N132").Value = -19092.8

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 70275
$ws.Range("J92").Value = 70275
$ws.Range("L92").Value = 70275
$ws.Range("N92").Value = -75267
$ws.Range("H117").Value = 88409
$ws.Range("J117").Value = 88409
$ws.Range("L117").Value = 88409
$ws.Range("N117").Value = -97587

Write-Host "Applied all cell updates."